$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: duplicate of row 2 (2025-03-24, 350, Food & Dining at, millennium)
$ws.Range("A2:D2").Copy($ws.Range("A4:D4"))

# Row 5: duplicate of row 2 again
$ws.Range("A2:D2").Copy($ws.Range("A5:D5"))

# Row 6: new date (2025-03-25) in A6, rest same as row 2 (B,C,D)
$ws.Range("B2:D2").Copy($ws.Range("B6:D6"))

# Stage the new date string as text (avoid Excel auto-converting it to a
# date serial number) in a scratch cell, then copy it into place as a
# shared string and strip the scratch formatting back off again.
$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = "2025-03-25"
$ws.Range("F1").Copy($ws.Range("A6"))
$ws.Range("A6").ClearFormats()
$ws.Range("F1").Clear()
